$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.155.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.009.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.697'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +12.00%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.744'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.99%  '
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000322'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.68'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.632.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.007.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.17%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.928.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '96.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.69%  '
$ws.Range("E28").Value = '  +19.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.131'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '13.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '48.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +19.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '673.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.445'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0819'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.81%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.150'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.93%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0487'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.149'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.65%  '
$ws.Range("E47").Value = '  -4.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '143.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000265'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.95%  '
